$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2").Value = 1.3
$ws.Range("P2").Value = 3.5
$ws.Range("Q2").Value = 1.97
$ws.Range("R2").Value = 1.93
